$d = $word.ActiveDocument
$d.Content.Find.Execute("869", $true, $false, $false, $false, $false, $true, 1, $false, "46", 2)
